$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "43.587.45"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.275.74"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.21%  "
Set-TextValue "D5" "122.51"
$ws.Range("E5").Value = "  +6.02%  "
Set-TextValue "D6" "266.28"
$ws.Range("E6").Value = "  -0.80%  "
Set-TextValue "D7" "0.643"
$ws.Range("E7").Value = "  +2.55%  "
$ws.Range("E8").Value = "  +0.18%  "
Set-TextValue "D9" "0.624"
$ws.Range("E9").Value = "  +0.79%  "
Set-TextValue "D10" "47.89"
$ws.Range("E10").Value = "  -2.16%  "
Set-TextValue "D11" "0.0948"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  -2.58%  "
Set-TextValue "D15" "0.910"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "2.617.09"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "2.260.78"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "43.546.17"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  -0.83%  "
Set-TextValue "D21" "72.22"
$ws.Range("E21").Value = "  -0.02%  "
Set-TextValue "D22" "2.41"
$ws.Range("E22").Value = "  -0.12%  "
Set-TextValue "D23" "235.36"
$ws.Range("E23").Value = "  +1.13%  "
Set-TextValue "D24" "9.53"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("E25").Value = "  -0.87%  "
Set-TextValue "D26" "12.01"
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("E27").Value = "  +1.71%  "
Set-TextValue "D28" "42.28"
$ws.Range("E28").Value = "  +0.43%  "
Set-TextValue "D29" "3.37"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +0.68%  "
Set-TextValue "D31" "172.42"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +0.45%  "
Set-TextValue "D33" "0.0916"
$ws.Range("E33").Value = "  -1.46%  "
Set-TextValue "D34" "5.73"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D35" "4.32"
$ws.Range("E35").Value = "  +13.97%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D36" "0.130"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("E37").Value = "  +4.85%  "
Set-TextValue "D38" "4.61"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("E40").Value = "  +5.13%  "
Set-TextValue "D41" "13.96"
$ws.Range("E41").Value = "  -4.13%  "
Set-TextValue "D42" "73.82"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -0.63%  "
Set-TextValue "D46" "5.70"
$ws.Range("E46").Value = "  -10.43%  "
Set-TextValue "D47" "74.00"
$ws.Range("E47").Value = "  +41.48%  "
$ws.Range("E48").Value = "  +0.04%  "
Set-TextValue "D49" "8.55"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("E50").Value = "  +0.29%  "
Set-TextValue "D51" "101.57"
$ws.Range("E51").Value = "  -1.22%  "
